# The workbook is already open; the active sheet ("vegetables") is the one
# that needs its G column dates normalized and its view updated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set every date in G2:G42 to the same serial date value as G1 (45011).
$ws.Range("G2:G42").Value = 45011

# Update the view: scroll so row 16 is at the top, and change the selection
# from G1:G42 to the single cell I36.
$win = $excel.ActiveWindow
$ws.Range("I36").Select()
$win.ScrollRow = 16
$win.ScrollColumn = 1
